$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -12.655
$ws.Range("E4").Value = 13.168
$ws.Range("E6").Value = 13.024
$ws.Range("A9").Value = -20.775
$ws.Range("E10").Value = 12.357
$ws.Range("C11").Value = -12.704
$ws.Range("E11").Value = 13.342
$ws.Range("A18").Value = -21.81
$ws.Range("A20").Value = -21.813
$ws.Range("D21").Value = -7.843999999999999
$ws.Range("E21").Value = 13.149
